$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.3663867549261397
$ws.Range("D2").Value = 0.01269013612205327
$ws.Range("E2").Value = 0.1289646226693755
$ws.Range("F2").Value = 8.714809173911618
$ws.Range("G2").Value = 0.002775424078496069
$ws.Range("J2").Value = 0.4657548065995769
$ws.Range("K2").Value = 6.130759347880144
$ws.Range("L2").Value = 0.09538396420233397
$ws.Range("C3").Value = 0.3651688647943985
$ws.Range("D3").Value = 0.01123078375112385
$ws.Range("E3").Value = 0.1295052471542029
$ws.Range("F3").Value = 8.586288342282359
$ws.Range("G3").Value = 0.002783461190924375
$ws.Range("J3").Value = 0.4629932560373575
$ws.Range("K3").Value = 6.057605493348262
$ws.Range("L3").Value = 0.09623760867909148
$ws.Range("C4").Value = 0.364619400146907
$ws.Range("D4").Value = 0.01033107699172575
$ws.Range("E4").Value = 0.1299005386054422
$ws.Range("F4").Value = 8.510250072220032
$ws.Range("G4").Value = 0.002788647399605783
$ws.Range("J4").Value = 0.4614707889218366
$ws.Range("K4").Value = 6.017315231203213
$ws.Range("L4").Value = 0.09680096102977132
$ws.Range("C5").Value = 0.3644452781432506
$ws.Range("D5").Value = 0.009963371461694237
$ws.Range("E5").Value = 0.130077567758871
$ws.Range("F5").Value = 8.479977656568025
$ws.Range("G5").Value = 0.002790824292718857
$ws.Range("J5").Value = 0.4608936753328123
$ws.Range("K5").Value = 6.002056789686378
$ws.Range("L5").Value = 0.09704042239521549
$ws.Range("C6").Value = 0.3644193701974103
$ws.Range("D6").Value = 0.009902245444571633
$ws.Range("E6").Value = 0.1301079268205783
$ws.Range("F6").Value = 8.474993798175348
$ws.Range("G6").Value = 0.002791189605120071
$ws.Range("J6").Value = 0.4608004548194629
$ws.Range("K6").Value = 5.99959311773739
$ws.Range("L6").Value = 0.09708078305528645
$ws.Range("C7").Value = 0.3646168503823901
$ws.Range("D7").Value = 0.0103261225044875
$ws.Range("E7").Value = 0.1299028615008169
$ws.Range("F7").Value = 8.509838929376542
$ws.Range("G7").Value = 0.002788676500512128
$ws.Range("J7").Value = 0.4614628307321738
$ws.Range("K7").Value = 6.01710475713395
$ws.Range("L7").Value = 0.09680415040083545
$ws.Range("C8").Value = 0.3659256124797139
$ws.Range("D8").Value = 0.01218762133152751
$ws.Range("E8").Value = 0.1291378915682593
$ws.Range("F8").Value = 8.669893248918811
$ws.Range("G8").Value = 0.002778143254595018
$ws.Range("J8").Value = 0.464766546875822
$ws.Range("K8").Value = 6.104573649419194
$ws.Range("L8").Value = 0.09567018090578472
$ws.Range("C9").Value = 0.3700702239161728
$ws.Range("D9").Value = 0.01581576951538466
$ws.Range("E9").Value = 0.1281398557129343
$ws.Range("F9").Value = 9.006988289050611
$ws.Range("G9").Value = 0.002759470472728474
$ws.Range("J9").Value = 0.4726302289338662
$ws.Range("K9").Value = 6.312980705007078
$ws.Range("L9").Value = 0.09375624352598244
$ws.Range("C10").Value = 0.3740848480057366
$ws.Range("D10").Value = 0.01847737615769773
$ws.Range("E10").Value = 0.1277121458084984
$ws.Range("F10").Value = 9.26943508674691
$ws.Range("G10").Value = 0.002746943850383044
$ws.Range("J10").Value = 0.4792690553316703
$ws.Range("K10").Value = 6.488858185524009
$ws.Range("L10").Value = 0.09253710574181184
$ws.Range("C11").Value = 0.376123534067915
$ws.Range("D11").Value = 0.01968943357982766
$ws.Range("E11").Value = 0.1275838409062473
$ws.Range("F11").Value = 9.392177357623041
$ws.Range("G11").Value = 0.002741500497276388
$ws.Range("J11").Value = 0.4824801033205546
$ws.Range("K11").Value = 6.573877054928744
$ws.Range("L11").Value = 0.09202272313312676
$ws.Range("C12").Value = 0.3769262144977006
$ws.Range("D12").Value = 0.02014876091337925
$ws.Range("E12").Value = 0.1275447773440774
$ws.Range("F12").Value = 9.439149862624731
$ws.Range("G12").Value = 0.00273947564362632
$ws.Range("J12").Value = 0.4837238129411361
$ws.Range("K12").Value = 6.60679716071229
$ws.Range("L12").Value = 0.09183369238776251
$ws.Range("C13").Value = 0.3767519765258385
$ws.Range("D13").Value = 0.02004981847520781
$ws.Range("E13").Value = 0.1275527669322187
$ws.Range("F13").Value = 9.429011414150352
$ws.Range("G13").Value = 0.002739910116032448
$ws.Range("J13").Value = 0.483454718728666
$ws.Range("K13").Value = 6.599674891846803
$ws.Range("L13").Value = 0.09187414802454086
$ws.Range("C14").Value = 0.3761889555468088
$ws.Range("D14").Value = 0.01972721480785822
$ws.Range("E14").Value = 0.1275804363020292
$ws.Range("F14").Value = 9.396031879873533
$ws.Range("G14").Value = 0.002741333182806172
$ws.Range("J14").Value = 0.4825818661231551
$ws.Range("K14").Value = 6.576570848055553
$ws.Range("L14").Value = 0.092007056280341
$ws.Range("C15").Value = 0.3758480875925443
$ws.Range("D15").Value = 0.01952966088766317
$ws.Range("E15").Value = 0.1275986245822232
$ws.Range("F15").Value = 9.375895434553684
$ws.Range("G15").Value = 0.002742209588202371
$ws.Range("J15").Value = 0.482050841949885
$ws.Range("K15").Value = 6.562513552943699
$ws.Range("L15").Value = 0.09208921506593271
$ws.Range("C16").Value = 0.3739559002942201
$ws.Range("D16").Value = 0.01839820503282397
$ws.Range("E16").Value = 0.1277218636836217
$ws.Range("F16").Value = 9.261481963125021
$ws.Range("G16").Value = 0.002747304697548235
$ws.Range("J16").Value = 0.4790630739164783
$ws.Range("K16").Value = 6.483403238443316
$ws.Range("L16").Value = 0.09257152862767626
$ws.Range("C17").Value = 0.3728495988606824
$ws.Range("D17").Value = 0.01770453692018492
$ws.Range("E17").Value = 0.1278144337682683
$ws.Range("F17").Value = 9.192159348448683
$ws.Range("G17").Value = 0.002750495531918418
$ws.Range("J17").Value = 0.4772793006117126
$ws.Range("K17").Value = 6.436158406201173
$ws.Range("L17").Value = 0.09287769139856827
$ws.Range("C18").Value = 0.3722332687468679
$ws.Range("D18").Value = 0.01730566788231869
$ws.Range("E18").Value = 0.1278739148305448
$ws.Range("F18").Value = 9.152601867397749
$ws.Range("G18").Value = 0.002752354838982589
$ws.Range("J18").Value = 0.4762712827371445
$ws.Range("K18").Value = 6.409455824370582
$ws.Range("L18").Value = 0.09305757441859797
$ws.Range("C19").Value = 0.3720280186025207
$ws.Range("D19").Value = 0.01717063170850253
$ws.Range("E19").Value = 0.1278951255376839
$ws.Range("F19").Value = 9.139262181367712
$ws.Range("G19").Value = 0.00275298850217032
$ws.Range("J19").Value = 0.4759330614247546
$ws.Range("K19").Value = 6.400495612468944
$ws.Range("L19").Value = 0.09311913085460333
$ws.Range("C20").Value = 0.3729652972265853
$ws.Range("D20").Value = 0.01777836645823072
$ws.Range("E20").Value = 0.1278039340446213
$ws.Range("F20").Value = 9.199506177692257
$ws.Range("G20").Value = 0.002750153377689139
$ws.Range("J20").Value = 0.4774673253299255
$ws.Range("K20").Value = 6.441138878320828
$ws.Range("L20").Value = 0.09284470818163371
$ws.Range("C21").Value = 0.3763534949007976
$ws.Range("D21").Value = 0.0198219605803942
$ws.Range("E21").Value = 0.1275720507324394
$ws.Range("F21").Value = 9.405705313876524
$ws.Range("G21").Value = 0.002740914207007464
$ws.Range("J21").Value = 0.4828374884860267
$ws.Range("K21").Value = 6.583337342783807
$ws.Range("L21").Value = 0.09196786193574624
$ws.Range("C22").Value = 0.3787467359318839
$ws.Range("D22").Value = 0.02115967023041776
$ws.Range("E22").Value = 0.1274760056036683
$ws.Range("F22").Value = 9.543344988300987
$ws.Range("G22").Value = 0.002735088077160164
$ws.Range("J22").Value = 0.486509119353471
$ws.Range("K22").Value = 6.680502222233315
$ws.Range("L22").Value = 0.09142832317813365
$ws.Range("C23").Value = 0.3774530080989962
$ws.Range("D23").Value = 0.0204454624699224
$ws.Range("E23").Value = 0.1275221896112768
$ws.Range("F23").Value = 9.469617412078208
$ws.Range("G23").Value = 0.00273817825919961
$ws.Range("J23").Value = 0.4845345882095131
$ws.Range("K23").Value = 6.628254851169913
$ws.Range("L23").Value = 0.09171322597366327
$ws.Range("C24").Value = 0.3729129286530508
$ws.Range("D24").Value = 0.01774498835493432
$ws.Range("E24").Value = 0.1278086614667266
$ws.Range("F24").Value = 9.196183751843563
$ws.Range("G24").Value = 0.002750307988109757
$ws.Range("J24").Value = 0.4773822648765815
$ws.Range("K24").Value = 6.438885777228052
$ws.Range("L24").Value = 0.09285960785024017
$ws.Range("C25").Value = 0.3687792932409479
$ws.Range("D25").Value = 0.01483565752631222
$ws.Range("E25").Value = 0.1283561602483694
$ws.Range("F25").Value = 8.91324462430876
$ws.Range("G25").Value = 0.002764311384139639
$ws.Range("J25").Value = 0.4703528773074765
$ws.Range("K25").Value = 6.25262479484752
$ws.Range("L25").Value = 0.09424104007939604
